$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 8940
$ws.Range("E2").Value = 787
$ws.Range("F2").Value = 787
$ws.Range("G2").Value = 670
$ws.Range("H2").Value = 379
$ws.Range("I2").Value = 238
$ws.Range("J2").Value = 141
$ws.Range("K2").Value = 10107
$ws.Range("L2").Value = 3020
$ws.Range("M2").Value = 7087
$ws.Range("N2").Value = 4304
$ws.Range("O2").Value = 2783
$ws.Range("P2").Value = 291
$ws.Range("Q2").Value = 466
$ws.Range("R2").Value = -1662
$ws.Range("S2").Value = 928
$ws.Range("T2").Value = 558
$ws.Range("U2").Value = -92
$ws.Range("V2").Value = 1673
$ws.Range("W2").Value = 8.8
$ws.Range("X2").Value = 4.24
$ws.Range("Y2").Value = 5.75
$ws.Range("Z2").Value = 3.97
$ws.Range("AA2").Value = 42.61
$ws.Range("AB2").Value = 1559.58
$ws.Range("AC2").Value = 410
$ws.Range("AD2").Value = 22.05
$ws.Range("AE2").Value = 10102
$ws.Range("AF2").Value = 0.89
$ws.Range("AG2").Value = 100
$ws.Range("AH2").Value = 1.11
$ws.Range("AI2").Value = 17.89
$ws.Range("AJ2").Value = 58141980

# Row 3
$ws.Range("D3").Value = 9968
$ws.Range("E3").Value = 797
$ws.Range("F3").Value = 797
$ws.Range("G3").Value = 803
$ws.Range("H3").Value = 546
$ws.Range("I3").Value = 386
$ws.Range("J3").Value = 160
$ws.Range("K3").Value = 12423
$ws.Range("L3").Value = 4138
$ws.Range("M3").Value = 8285
$ws.Range("N3").Value = 4654
$ws.Range("O3").Value = 3630
$ws.Range("P3").Value = 291
$ws.Range("Q3").Value = 749
$ws.Range("R3").Value = -1175
$ws.Range("S3").Value = 631
$ws.Range("T3").Value = 1165
$ws.Range("U3").Value = -416
$ws.Range("V3").Value = 2410
$ws.Range("W3").Value = 7.99
$ws.Range("X3").Value = 5.48
$ws.Range("Y3").Value = 8.62
$ws.Range("Z3").Value = 4.85
$ws.Range("AA3").Value = 49.95
$ws.Range("AB3").Value = 1683.31
$ws.Range("AC3").Value = 664
$ws.Range("AD3").Value = 17.04
$ws.Range("AE3").Value = 10925
$ws.Range("AF3").Value = 1.04
$ws.Range("AG3").Value = 100
$ws.Range("AH3").Value = 0.88
$ws.Range("AI3").Value = 11.03
$ws.Range("AJ3").Value = 58141980

# Row 4
$ws.Range("D4").Value = 10150
$ws.Range("E4").Value = 610
$ws.Range("F4").Value = 610
$ws.Range("G4").Value = 572
$ws.Range("H4").Value = 346
$ws.Range("I4").Value = 198
$ws.Range("J4").Value = 147
$ws.Range("K4").Value = 14673
$ws.Range("L4").Value = 6027
$ws.Range("M4").Value = 8646
$ws.Range("N4").Value = 4858
$ws.Range("O4").Value = 3788
$ws.Range("P4").Value = 291
$ws.Range("Q4").Value = 188
$ws.Range("R4").Value = -1818
$ws.Range("S4").Value = 1687
$ws.Range("T4").Value = 1490
$ws.Range("U4").Value = -1302
$ws.Range("V4").Value = 4131
$ws.Range("W4").Value = 6.01
$ws.Range("X4").Value = 3.41
$ws.Range("Y4").Value = 4.17
$ws.Range("Z4").Value = 2.55
$ws.Range("AA4").Value = 69.7
$ws.Range("AB4").Value = 1750.82
$ws.Range("AC4").Value = 341
$ws.Range("AD4").Value = 25.83
$ws.Range("AE4").Value = 11250
$ws.Range("AF4").Value = 0.78
$ws.Range("AG4").Value = 100
$ws.Range("AH4").Value = 1.14
$ws.Range("AI4").Value = 21.78
$ws.Range("AJ4").Value = 58141980

# Row 5
$ws.Range("D5").Value = 11240
$ws.Range("E5").Value = 845
$ws.Range("F5").Value = 845
$ws.Range("G5").Value = 835
$ws.Range("H5").Value = 857
$ws.Range("I5").Value = 663
$ws.Range("J5").Value = 195
$ws.Range("K5").Value = 15492
$ws.Range("L5").Value = 6098
$ws.Range("M5").Value = 9393
$ws.Range("N5").Value = 5446
$ws.Range("O5").Value = 3947
$ws.Range("P5").Value = 291
$ws.Range("Q5").Value = 1478
$ws.Range("R5").Value = -825
$ws.Range("S5").Value = -348
$ws.Range("T5").Value = 746
$ws.Range("U5").Value = 732
$ws.Range("V5").Value = 3851
$ws.Range("W5").Value = 7.52
$ws.Range("X5").Value = 7.63
$ws.Range("Y5").Value = 12.86
$ws.Range("Z5").Value = 5.68
$ws.Range("AA5").Value = 64.92
$ws.Range("AB5").Value = 1960.64
$ws.Range("AC5").Value = 1139
$ws.Range("AD5").Value = 16.98
$ws.Range("AE5").Value = 12611
$ws.Range("AF5").Value = 1.53
$ws.Range("AG5").Value = 100
$ws.Range("AH5").Value = 0.52
$ws.Range("AI5").Value = 6.52
$ws.Range("AJ5").Value = 58141980

# Row 6
$ws.Range("D6").Value = 12313
$ws.Range("E6").Value = 858
$ws.Range("F6").Value = 858
$ws.Range("G6").Value = 608
$ws.Range("H6").Value = 327
$ws.Range("I6").Value = 299
$ws.Range("K6").Value = 15922
$ws.Range("L6").Value = 6279
$ws.Range("M6").Value = 9643
$ws.Range("N6").Value = 5696
$ws.Range("P6").Value = 291
$ws.Range("Q6").Value = 454
$ws.Range("R6").Value = -792
$ws.Range("S6").Value = -64
$ws.Range("T6").Value = 657
$ws.Range("U6").Value = -203
$ws.Range("V6").Value = 3868
$ws.Range("W6").Value = 6.97
$ws.Range("X6").Value = 2.66
$ws.Range("Y6").Value = 5.36
$ws.Range("Z6").Value = 2.08
$ws.Range("AA6").Value = 65.11
$ws.Range("AB6").Value = 2056.69
$ws.Range("AC6").Value = 514
$ws.Range("AD6").Value = 34.66
$ws.Range("AE6").Value = 13191
$ws.Range("AF6").Value = 1.35
$ws.Range("AG6").Value = 100
$ws.Range("AH6").Value = 0.56
$ws.Range("AI6").Value = 14.46
$ws.Range("AJ6").Value = 58141980

# Rows 7-9: clear all data columns (D through AJ), keep only A/B/C
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()